$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.821.36"
$ws.Range("E2").Value = "  +4.12%  "
$ws.Range("D3").Value = "1.877.21"
$ws.Range("E3").Value = "  +3.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "278.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5259"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3452"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06952"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8040"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07867"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "1.926.56"
$ws.Range("E14").Value = "  +6.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.158"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.63%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.01%  "
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008075"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "26.856.58"
$ws.Range("E21").Value = "  +4.11%  "
$ws.Range("D22").Value = "2.124.55"
$ws.Range("E22").Value = "  +4.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.749"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.191"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.345"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.662"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.358"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.327"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08909"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.175"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7347"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.903"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.249"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.395"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01850"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5132"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9617"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.221"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.074"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4498"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1344"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.333"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05944"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.06%  "
